# Datos_Reservas.xlsx - "mesas" (tables) seed data rewrite
#
# Before: 4 rows of table data with non-contiguous/inconsistent IDs
#   (2, 3, 5, 10) and a stray text-typed Capacidad value on the last
#   row, plus a highlighted/bordered header row.
# After : 9 contiguous rows (IDs 1-9), all "Libre" / Capacidad=4,
#   numeric Capacidad throughout, and a plain (unstyled) header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row: drop the fill/border highlight style -------------------
$ws.Range("A1:D1").ClearFormats()

# ---- Table rows: ID, Usuario, Estado, Capacidad -------------------------
# Usuario ("B") is always blank (an empty, text-typed placeholder cell -
# matches the "free table" rows before anyone reserves them).
$estado = "Libre"
$capacidad = 4

for ($i = 1; $i -le 9; $i++) {
    $row = $i + 1

    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Formula = "'"
    $ws.Cells.Item($row, 3).Value = $estado
    $ws.Cells.Item($row, 4).Value = $capacidad
}

# ---- View: plain single-cell selection on A1 -----------------------------
$ws.Range("A1").Select()
